$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4: the "age" column (C) switches from a text value ("30") to a real
# number (30). Re-point the cell format to a general, right-aligned,
# wrap-text number format (matching the sibling numeric column F) before
# writing the numeric value, so it is stored as t="n" instead of t="s".
foreach ($r in 2..4) {
    $c = $ws.Cells.Item($r, 3)
    $c.NumberFormat = "General"
    $c.HorizontalAlignment = -4152
    $c.WrapText = $true
    $c.Value = 30
}

# Rows 2-4: the "create time" (M) and "create time's time part" (K) columns
# were regenerated a few minutes later.
foreach ($r in 2..4) {
    $ws.Cells.Item($r, 11).Value = "11:12:10"
    $ws.Cells.Item($r, 13).Value = "2019-10-11 11:12:10"
}
